$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.939725160598755
$ws.Range("B1").Value = 4.608788967132568
$ws.Range("C1").Value = 3.346324920654297
$ws.Range("D1").Value = 2.519932508468628
$ws.Range("E1").Value = 2.030664443969727
